$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the current row 114 (Puerro,
# Vega Central Mapocho de Santiago). Shift the existing rows 114:124 down
# to 115:125 and fill the vacated row 114 with the new record.
$ws.Rows("114:114").Insert()

$ws.Range("A114").Value2 = 9
$ws.Range("B114").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C114").Value2 = "Metropolitana"
$ws.Range("D114").Value2 = 45077
$ws.Range("E114").Value2 = 13
$ws.Range("F114").Value2 = 100112005
$ws.Range("G114").Value2 = "Puerro"
$ws.Range("H114").Value2 = "Sin especificar"
$ws.Range("I114").Value2 = "Primera"
$ws.Range("J114").Value2 = 70
$ws.Range("K114").Value2 = 8000
$ws.Range("L114").Value2 = 8000
$ws.Range("M114").Value2 = 8000
$ws.Range("N114").Value2 = "$/paquete 20 unidades"
$ws.Range("O114").Value2 = "Provincia de Chacabuco"
$ws.Range("P114").Value2 = 400
$ws.Range("Q114").Value2 = 20
$ws.Range("R114").Value2 = "Hortaliza"
